$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New bold cell D8: ratio formula, formatted like the other bold cells (A8/B8/C8 row)
$ws.Range("D8").Formula = "=PI()*B2^2/4/B3/SIN(C1)"
$ws.Range("D8").Font.Bold = $true

# New row 13: label "S" and its surface-area formula
$ws.Range("A13").Value = "S"
$ws.Range("B13").Formula = "=PI()*(B2/1000)^2/4"

# Update the active selection shown in the saved view
$ws.Range("O12").Select()
